# Apply the commit "updated with tests, new code (lstm)" changes:
#  1) Rename several header labels in row 1
#  2) Remove the stray ".." placeholder value from C79 (was a text marker)
#     by replacing it with a real numeric GDP estimate
#  3) Update the "GDP" column (C) data rows with refreshed (imputed) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) data updates ---
$ws.Range("C3").Value = 5596.139681459835
$ws.Range("C4").Value = 2934.187009790061
$ws.Range("C5").Value = 2870.311589353206
$ws.Range("C6").Value = 1873.394108966653
$ws.Range("C7").Value = 8947.741473873051
$ws.Range("C8").Value = 1460.056109840828
$ws.Range("C9").Value = 1909.084588129339
$ws.Range("C10").Value = 10594.98659239237
$ws.Range("C11").Value = 6128.19547247793
$ws.Range("C13").Value = 471.181692645893
$ws.Range("C15").Value = 2898.942214704482
$ws.Range("C16").Value = 1904.346464968814
$ws.Range("C17").Value = 5730.354774594881
$ws.Range("C18").Value = 1503.870423231357
$ws.Range("C19").Value = 9271.398233246389
$ws.Range("C20").Value = 2983.242707849043
$ws.Range("C21").Value = 1955.461557360978
$ws.Range("C22").Value = 11286.24301624575
$ws.Range("C23").Value = 6336.709213679884
$ws.Range("C25").Value = 492.3430015592067
$ws.Range("C27").Value = 2965.153206179127
$ws.Range("C28").Value = 1939.33862702996
$ws.Range("C29").Value = 1577.487171555845
$ws.Range("C30").Value = 3083.80337578809
$ws.Range("C31").Value = 5885.254624554112
$ws.Range("C32").Value = 9477.887185090232
$ws.Range("C33").Value = 2024.117324382548
$ws.Range("C34").Value = 11627.81065059172
$ws.Range("C35").Value = 6711.616186806423
$ws.Range("C37").Value = 2999.422762626143
$ws.Range("C38").Value = 3156.723844635973
$ws.Range("C39").Value = 1657.651524528445
$ws.Range("C40").Value = 1982.009737844954
$ws.Range("C41").Value = 2094.024217383061
$ws.Range("C42").Value = 6911.59200404802
$ws.Range("C43").Value = 6051.685746144485
$ws.Range("C44").Value = 9690.869064532331
$ws.Range("C47").Value = 513.7390871590731
$ws.Range("C48").Value = 3056.152683606517
$ws.Range("C49").Value = 3212.740625904757
$ws.Range("C50").Value = 1716.389195271215
$ws.Range("C51").Value = 2000.792448761861
$ws.Range("C52").Value = 2201.396847776877
$ws.Range("C53").Value = 7200.731056811853
$ws.Range("C54").Value = 6203.843262938323
$ws.Range("C55").Value = 9693.722968944676
$ws.Range("C58").Value = 534.5063430177229
$ws.Range("C59").Value = 2286.013198234259
$ws.Range("C60").Value = 7449.08671983612
$ws.Range("C61").Value = 6255.426161047989
$ws.Range("C62").Value = 3252.634165082374
$ws.Range("C63").Value = 3137.260298393558
$ws.Range("C64").Value = 2025.814194788851
$ws.Range("C65").Value = 558.2093442539386
$ws.Range("C68").Value = 2361.056581219794
$ws.Range("C69").Value = 7580.275568826287
$ws.Range("C70").Value = 6522.736799041846
$ws.Range("C71").Value = 3314.741082534716
$ws.Range("C72").Value = 3210.869677115934
$ws.Range("C73").Value = 2067.29003376698
$ws.Range("C74").Value = 579.0880693780265
$ws.Range("C77").Value = 2425.561644739583
$ws.Range("C78").Value = 7633.969039669125
$ws.Range("C79").Value = 6550.274372976741
$ws.Range("C80").Value = 3382.563653843273
$ws.Range("C81").Value = 3242.636921959078
$ws.Range("C82").Value = 2111.193164269742
$ws.Range("C83").Value = 1895.214690888655
